$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Add "PA_Faltante_COT" calculated column (H) ---
$colCOT = $lo.ListColumns.Add()
$ws.Cells.Item(1, 8).Value = "PA_Faltante_COT"
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=(5431200-Tabla1354[[#This Row],[PA_Acumulada]])"
}

# --- Add "PA_Faltante_TOT" calculated column (I) ---
$colTOT = $lo.ListColumns.Add()
$ws.Cells.Item(1, 9).Value = "PA_Faltante_TOT"
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=(10862400-Tabla1354[[#This Row],[PA_Acumulada]])"
}

# --- Match the look of the new columns to the existing calculated column
#     (PA_Faltante_Miembro, col G): header style + the two body-row styles
#     (regular rows vs. the bold "MDRT" rows at the bottom of the table). ---
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

$ws.Range("G2:G42").Copy()
$ws.Range("H2:I42").PasteSpecial(-4122)

$ws.Range("G43:G51").Copy()
$ws.Range("H43:I51").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Restore a plain view: clear the old scroll position and move the
#     selection to the first new calculated cell, like after typing it in. ---
$ws.Range("I3").Select()
